$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (Strike#) values per row, regenerated for this commit
$kValues = @{
    2 = 0
    3 = 1
    4 = 0
    5 = 0
    6 = 1
    7 = 0
    8 = 1
    9 = 1
    10 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 2
    22 = 0
    23 = 0
    24 = 0
    25 = 2
    26 = 0
    27 = 2
    28 = 1
    29 = 1
    30 = 0
    31 = 0
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 0
    42 = 1
    43 = 1
    44 = 0
    45 = 3
    46 = 1
    47 = 1
    48 = 0
    49 = 3
    50 = 1
    51 = 1
    52 = 1
    53 = 1
    54 = 1
    55 = 1
    56 = 2
    57 = 1
    58 = 0
    59 = 2
    60 = 1
    61 = 1
    62 = 0
    63 = 0
    64 = 2
    65 = 1
    66 = 1
    67 = 1
    68 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
